# config change for installed capacity
# Updates the literal "installed capacity" numbers feeding the
# regional_profile / ists_gen / state_gen / max_info report sheets.
# Dependent SUM/addition formulas on the same sheets recompute automatically.

$wb = $excel.ActiveWorkbook

# --- regional_profile (sheet2) ---------------------------------------
$wsRegional = $wb.Worksheets.Item("regional_profile")
$wsRegional.Range("B9").Value = 15614
$wsRegional.Range("B10").Value = 6382

# --- ists_gen (sheet3) -------------------------------------------------
$wsIsts = $wb.Worksheets.Item("ists_gen")
$wsIsts.Range("B16").Value = 243

# --- state_gen (sheet4) -------------------------------------------------
$wsState = $wb.Worksheets.Item("state_gen")
$wsState.Range("B6").Value = 2492
$wsState.Range("B7").Value = 6299
$wsState.Range("B10").Value = 1580
$wsState.Range("B11").Value = 5010
$wsState.Range("B14").Value = 1560
$wsState.Range("B15").Value = 2424

# --- max_info (sheet5) - mirrors the same figures as literal values ----
$wsMax = $wb.Worksheets.Item("max_info")
$wsMax.Range("B16").Value = 243
$wsMax.Range("B29").Value = 2492
$wsMax.Range("B30").Value = 6299
$wsMax.Range("B33").Value = 1580
$wsMax.Range("B34").Value = 5010
$wsMax.Range("B37").Value = 1560
$wsMax.Range("B38").Value = 2424
$wsMax.Range("B45").Value = 15614
$wsMax.Range("B46").Value = 6382

# Recalculate so every dependent formula's cached <v> reflects the new inputs
$excel.Calculate()

# --- Selection / active-sheet bookkeeping ------------------------------
# Touch ists_gen, state_gen and max_info first (each becomes active only
# transiently) so their scroll/selection state updates without leaving
# them as the final active tab.
[void]$wsIsts.Activate()
$excel.ActiveWindow.Zoom = 130
[void]$wsIsts.Range("B21").Select()

[void]$wsState.Activate()
[void]$wsState.Range("B19").Select()

[void]$wsMax.Activate()
[void]$wsMax.Range("B47").Select()

# regional_profile is activated last, matching the workbook's final
# activeTab / tabSelected state.
[void]$wsRegional.Activate()
$excel.ActiveWindow.Zoom = 130
[void]$wsRegional.Range("B11").Select()
